$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 211.4614666666667
$ws.Range("H2").Value = 634.3844
$ws.Range("I2").Value = 0.2421062275331183
$ws.Range("J2").Value = 0.2421062275331183
$ws.Range("M2").Value = 0.03992133333333333
$ws.Range("O2").Value = 0.4073688578678476
$ws.Range("P2").Value = 0.4073688578678477
$ws.Range("Q2").Value = 8.441823697955554
$ws.Range("R2").Value = 75.9764132816
$ws.Range("S2").Value = 0.09862653739285965
$ws.Range("T2").Value = 0.09862653739285966
$ws.Range("G3").Value = 211.4614666666667
$ws.Range("H3").Value = 634.3844
$ws.Range("I3").Value = 0.2421062275331183
$ws.Range("J3").Value = 0.2421062275331183
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05807666666666667
$ws.Range("N3").Value = 0.17423
$ws.Range("O3").Value = 0.5926311421321523
$ws.Range("P3").Value = 0.5926311421321524
$ws.Range("Q3").Value = 12.28097711244444
$ws.Range("R3").Value = 110.528794012
$ws.Range("S3").Value = 0.1434796901402586
$ws.Range("T3").Value = 0.1434796901402587
$ws.Range("I4").Value = 0.08842543241393927
$ws.Range("J4").Value = 0.08842543241393927
$ws.Range("M4").Value = 0.03992133333333333
$ws.Range("O4").Value = 0.4073688578678476
$ws.Range("P4").Value = 0.4073688578678477
$ws.Range("Q4").Value = 3.083241263390666
$ws.Range("S4").Value = 0.03602176740893699
$ws.Range("T4").Value = 0.036021767408937
$ws.Range("I5").Value = 0.08842543241393927
$ws.Range("J5").Value = 0.08842543241393927
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05807666666666667
$ws.Range("N5").Value = 0.17423
$ws.Range("O5").Value = 0.5926311421321523
$ws.Range("P5").Value = 0.5926311421321524
$ws.Range("Q5").Value = 4.485430724763333
$ws.Range("R5").Value = 40.36887652287
$ws.Range("S5").Value = 0.05240366500500227
$ws.Range("T5").Value = 0.05240366500500228
$ws.Range("G6").Value = 174.3107043333333
$ws.Range("H6").Value = 522.932113
$ws.Range("I6").Value = 0.199571617988009
$ws.Range("J6").Value = 0.199571617988009
$ws.Range("M6").Value = 0.03992133333333333
$ws.Range("O6").Value = 0.4073688578678476
$ws.Range("P6").Value = 0.4073688578678477
$ws.Range("Q6").Value = 6.95871573125911
$ws.Range("R6").Value = 62.62844158133199
$ws.Range("S6").Value = 0.08129926208261364
$ws.Range("T6").Value = 0.08129926208261364
$ws.Range("G7").Value = 174.3107043333333
$ws.Range("H7").Value = 522.932113
$ws.Range("I7").Value = 0.199571617988009
$ws.Range("J7").Value = 0.199571617988009
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.05807666666666667
$ws.Range("N7").Value = 0.17423
$ws.Range("O7").Value = 0.5926311421321523
$ws.Range("P7").Value = 0.5926311421321524
$ws.Range("Q7").Value = 10.12338467199889
$ws.Range("R7").Value = 91.11046204798998
$ws.Range("S7").Value = 0.1182723559053954
$ws.Range("T7").Value = 0.1182723559053954
$ws.Range("G8").Value = 28.53474833333333
$ws.Range("H8").Value = 85.60424499999999
$ws.Range("I8").Value = 0.03266997236655063
$ws.Range("J8").Value = 0.03266997236655063
$ws.Range("M8").Value = 0.03992133333333333
$ws.Range("O8").Value = 0.4073688578678476
$ws.Range("P8").Value = 0.4073688578678477
$ws.Range("Q8").Value = 1.139145199797778
$ws.Range("R8").Value = 10.25230679818
$ws.Range("S8").Value = 0.01330872932953587
$ws.Range("T8").Value = 0.01330872932953587
$ws.Range("G9").Value = 28.53474833333333
$ws.Range("H9").Value = 85.60424499999999
$ws.Range("I9").Value = 0.03266997236655063
$ws.Range("J9").Value = 0.03266997236655063
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05807666666666667
$ws.Range("N9").Value = 0.17423
$ws.Range("O9").Value = 0.5926311421321523
$ws.Range("P9").Value = 0.5926311421321524
$ws.Range("Q9").Value = 1.657203067372222
$ws.Range("R9").Value = 14.91482760635
$ws.Range("S9").Value = 0.01936124303701475
$ws.Range("T9").Value = 0.01936124303701476
$ws.Range("G10").Value = 230.32901
$ws.Range("H10").Value = 690.98703
$ws.Range("I10").Value = 0.263708034289011
$ws.Range("J10").Value = 0.263708034289011
$ws.Range("M10").Value = 0.03992133333333333
$ws.Range("O10").Value = 0.4073688578678476
$ws.Range("P10").Value = 0.4073688578678477
$ws.Range("Q10").Value = 9.195041184546666
$ws.Range("R10").Value = 82.75537066091999
$ws.Range("S10").Value = 0.1074264407388896
$ws.Range("T10").Value = 0.1074264407388896
$ws.Range("G11").Value = 230.32901
$ws.Range("H11").Value = 690.98703
$ws.Range("I11").Value = 0.263708034289011
$ws.Range("J11").Value = 0.263708034289011
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.05807666666666667
$ws.Range("N11").Value = 0.17423
$ws.Range("O11").Value = 0.5926311421321523
$ws.Range("P11").Value = 0.5926311421321524
$ws.Range("Q11").Value = 13.37674113743333
$ws.Range("R11").Value = 120.3906702369
$ws.Range("S11").Value = 0.1562815935501213
$ws.Range("T11").Value = 0.1562815935501214
$ws.Range("G12").Value = 151.5554656666667
$ws.Range("H12").Value = 454.666397
$ws.Range("I12").Value = 0.1735187154093718
$ws.Range("J12").Value = 0.1735187154093718
$ws.Range("M12").Value = 0.03992133333333333
$ws.Range("O12").Value = 0.4073688578678476
$ws.Range("P12").Value = 0.4073688578678477
$ws.Range("Q12").Value = 6.050296263367555
$ws.Range("R12").Value = 54.45266637030799
$ws.Range("S12").Value = 0.07068612091501188
$ws.Range("T12").Value = 0.07068612091501189
$ws.Range("G13").Value = 151.5554656666667
$ws.Range("H13").Value = 454.666397
$ws.Range("I13").Value = 0.1735187154093718
$ws.Range("J13").Value = 0.1735187154093718
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.05807666666666667
$ws.Range("N13").Value = 0.17423
$ws.Range("O13").Value = 0.5926311421321523
$ws.Range("P13").Value = 0.5926311421321524
$ws.Range("Q13").Value = 8.801836261034444
$ws.Range("R13").Value = 79.21652634930999
$ws.Range("S13").Value = 0.1028325944943599
$ws.Range("T13").Value = 0.1028325944943599
